$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.196.57"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "2.826.58"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'356.17"
$ws.Range("E5").Value = "  +2.60%  "
$ws.Range("D6").Value = "'112.35"
$ws.Range("E6").Value = "  -3.69%  "
$ws.Range("D7").Value = "'0.572"
$ws.Range("E7").Value = "  +3.29%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.600"
$ws.Range("E9").Value = "  +2.16%  "
$ws.Range("D10").Value = "'40.92"
$ws.Range("E10").Value = "  -5.39%  "
$ws.Range("D11").Value = "'0.0863"
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "'19.99"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "'7.77"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").Value = "3.265.23"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "2.833.40"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("E17").Value = "  +4.38%  "
$ws.Range("D18").Value = "51.970.95"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "'7.56"
$ws.Range("E19").Value = "  +6.18%  "
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").Value = "0.0₃0997"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").Value = "'70.91"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").Value = "'270.67"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "'2.84"
$ws.Range("E25").Value = "  +2.71%  "
$ws.Range("D26").Value = "'27.02"
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "'10.34"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").Value = "'0.145"
$ws.Range("E30").Value = "  +3.07%  "
$ws.Range("D31").Value = "'0.0489"
$ws.Range("E31").Value = "  +19.05%  "
$ws.Range("D32").Value = "'52.43"
$ws.Range("E32").Value = "  +4.25%  "
$ws.Range("D33").Value = "'34.82"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").Value = "'5.96"
$ws.Range("E34").Value = "  +4.15%  "
$ws.Range("D35").Value = "'5.65"
$ws.Range("E35").Value = "  +13.34%  "
$ws.Range("D36").Value = "'0.0856"
$ws.Range("E36").Value = "  +3.75%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("D39").Value = "'2.05"
$ws.Range("E39").Value = "  -3.69%  "
$ws.Range("D40").Value = "'18.35"
$ws.Range("E40").Value = "  -3.13%  "
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("D42").Value = "'127.18"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").Value = "'23.26"
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("D44").Value = "'2.51"
$ws.Range("E44").Value = "  -7.54%  "
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("D46").Value = "'3.37"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").Value = "2.086.95"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("E48").Value = "  -5.56%  "
$ws.Range("D49").Value = "'5.90"
$ws.Range("E49").Value = "  +6.50%  "
$ws.Range("D50").Value = "'0.977"
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("D51").Value = "'9.16"
$ws.Range("E51").Value = "  +2.22%  "
